# Portfolio Investment data preparation
# - insert a new "shares" column between Buy-in-date and Buy-in-amount
# - rename headers to lowercase/snake_case and add an "exit_date" header
# - fill in the new shares values for the three existing holdings

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("shares"); this shifts the existing "Buy-in-amount"
# column (and its data / formatting) from D to E without touching any of the
# already-entered text (equity/ticker/date strings keep their original
# shared-string entries and styling).
$ws.Columns.Item(4).Insert()

# New "shares" values for the three holdings.
$ws.Range("D2").Value = 12
$ws.Range("D3").Value = 20
$ws.Range("D4").Value = 100

# Re-label the header row (lowercase / snake_case), including the new
# "shares" column and the (previously-unused) "exit_date" column.
$ws.Range("A1").Value = "equity"
$ws.Range("B1").Value = "ticker"
$ws.Range("C1").Value = "buy_in_date"
$ws.Range("D1").Value = "shares"
$ws.Range("E1").Value = "buy_in_amount"
$ws.Range("F1").Value = "exit_date"

# Reflect the resulting selection location.
$ws.Range("B13").Select()

$wb.Save()
